function Set-TextValue($ws, $r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.197.93'
$ws.Cells.Item(2, 5).Value = '  +0.32%  '

$ws.Cells.Item(3, 4).Value = '1.843.37'
$ws.Cells.Item(3, 5).Value = '  +0.61%  '

Set-TextValue $ws 4 4 '0.9993'
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

Set-TextValue $ws 5 4 '240.71'
$ws.Cells.Item(5, 5).Value = '  +0.06%  '

Set-TextValue $ws 6 4 '0.6706'
$ws.Cells.Item(6, 5).Value = '  -2.03%  '

$ws.Cells.Item(7, 5).Value = '  -0.01%  '

Set-TextValue $ws 8 4 '0.07420'
$ws.Cells.Item(8, 5).Value = '  -0.56%  '

$ws.Cells.Item(9, 5).Value = '  -1.98%  '

Set-TextValue $ws 10 4 '22.84'
$ws.Cells.Item(10, 5).Value = '  -1.09%  '

$ws.Cells.Item(12, 4).Value = '1.823.30'
$ws.Cells.Item(12, 5).Value = '  -0.43%  '

Set-TextValue $ws 13 4 '5.003'
$ws.Cells.Item(13, 5).Value = '  -1.12%  '

Set-TextValue $ws 14 4 '0.6745'
$ws.Cells.Item(14, 5).Value = '  -1.09%  '

Set-TextValue $ws 15 4 '86.12'
$ws.Cells.Item(15, 5).Value = '  -1.80%  '

$ws.Cells.Item(16, 5).Value = '  -0.17%  '

$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws 17 4 '0.000008309'
$ws.Cells.Item(17, 5).Value = '  +1.56%  '

$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '29.111.70'
$ws.Cells.Item(18, 5).Value = '  +0.05%  '

$ws.Cells.Item(19, 5).Value = '  +0.10%  '

Set-TextValue $ws 20 4 '12.52'
$ws.Cells.Item(20, 5).Value = '  -0.02%  '

Set-TextValue $ws 21 4 '1.000'
$ws.Cells.Item(21, 5).Value = '  +0.01%  '

Set-TextValue $ws 22 4 '7.182'
$ws.Cells.Item(22, 5).Value = '  -3.31%  '

$ws.Cells.Item(23, 5).Value = '  +0.01%  '

Set-TextValue $ws 24 4 '160.44'
$ws.Cells.Item(24, 5).Value = '  +0.16%  '

Set-TextValue $ws 25 4 '8.688'
$ws.Cells.Item(25, 5).Value = '  -0.67%  '

Set-TextValue $ws 26 4 '0.1402'
$ws.Cells.Item(26, 5).Value = '  -3.86%  '

$ws.Cells.Item(27, 5).Value = '  -0.68%  '

Set-TextValue $ws 28 4 '1.507'
$ws.Cells.Item(28, 5).Value = '  -0.11%  '

Set-TextValue $ws 29 4 '4.177'
$ws.Cells.Item(29, 5).Value = '  -2.70%  '

Set-TextValue $ws 30 4 '4.067'
$ws.Cells.Item(30, 5).Value = '  -1.97%  '

Set-TextValue $ws 32 4 '0.05313'
$ws.Cells.Item(32, 5).Value = '  +2.90%  '

Set-TextValue $ws 33 4 '0.7578'
$ws.Cells.Item(33, 5).Value = '  -1.16%  '

$ws.Cells.Item(34, 5).Value = '  +1.68%  '

Set-TextValue $ws 35 4 '1.136'
$ws.Cells.Item(35, 5).Value = '  +0.27%  '

Set-TextValue $ws 36 4 '2.678'
$ws.Cells.Item(36, 5).Value = '  +0.10%  '

$ws.Cells.Item(37, 4).Value = '1.326.64'
$ws.Cells.Item(37, 5).Value = '  +1.49%  '

Set-TextValue $ws 38 4 '0.01804'
$ws.Cells.Item(38, 5).Value = '  -1.88%  '

Set-TextValue $ws 39 4 '2.730'
$ws.Cells.Item(39, 5).Value = '  +0.41%  '

Set-TextValue $ws 40 4 '0.9226'
$ws.Cells.Item(40, 5).Value = '  -0.77%  '

Set-TextValue $ws 41 4 '5.948'
$ws.Cells.Item(41, 5).Value = '  +2.02%  '

$ws.Cells.Item(42, 5).Value = '  +0.12%  '

Set-TextValue $ws 43 4 '103.42'
$ws.Cells.Item(43, 5).Value = '  -1.15%  '

Set-TextValue $ws 44 4 '0.07912'
$ws.Cells.Item(44, 5).Value = '  +13.07%  '

$ws.Cells.Item(45, 4).Value = '1.967.52'
$ws.Cells.Item(45, 5).Value = '  -0.97%  '

Set-TextValue $ws 46 4 '0.5163'
$ws.Cells.Item(46, 5).Value = '  -0.65%  '

$ws.Cells.Item(47, 5).Value = '  +0.00%  '

$ws.Cells.Item(48, 5).Value = '  -1.86%  '

Set-TextValue $ws 49 4 '63.69'
$ws.Cells.Item(49, 5).Value = '  -2.35%  '

Set-TextValue $ws 50 4 '9.167'
$ws.Cells.Item(50, 5).Value = '  -4.01%  '

Set-TextValue $ws 51 4 '0.05945'
$ws.Cells.Item(51, 5).Value = '  +0.41%  '

